# Fix typo 'Cyro-EM' -> 'Cryo-EM' in the storage_medium lookup sheet,
# and refresh the pav:createdOn timestamp on the .metadata sheet.
# (Closes #9)

$wb = $excel.ActiveWorkbook

# --- 1. storage_medium sheet: fix the typo and re-sequence the
#        Cryo-EM / DMSO (serum) / RNAlater rows (12-14) ---
$storage = $wb.Worksheets.Item("storage_medium")

$storage.Range("A12").Value = "Cryo-EM"
$storage.Range("B12").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333"

$storage.Range("A13").Value = "DMSO (serum)"
$storage.Range("B13").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125"

$storage.Range("A14").Value = "RNAlater"
$storage.Range("B14").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"

# --- 2. .metadata sheet: update pav:createdOn timestamp ---
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2024-03-14T10:53:58-04:00"
